$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Part 1: update rows 2-41 (price/volume refresh) ----
$dTextRows = @(2, 3, 5, 6, 7, 8, 9, 10, 12, 13, 14, 15, 17, 18, 19, 21, 22, 23, 24, 25, 26, 27, 28, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41)
foreach ($r in $dTextRows) { $ws.Cells.Item($r, 4).NumberFormat = "@" }

$ws.Range("D2").Value = "27.797.36"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "1.879.59"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "332.21"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.4718"
$ws.Range("E7").Value = "  +4.30%  "
$ws.Range("D8").Value = "0.3955"
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("D9").Value = "47.87"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").Value = "0.08064"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "22.28"
$ws.Range("E12").Value = "  +4.37%  "
$ws.Range("D13").Value = "1.876.79"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "5.978"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "7.143"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "0.00001050"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "87.32"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "0.06667"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "27.794.21"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").Value = "5.538"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "11.00"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "2.305"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "2.108.97"
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("D27").Value = "159.51"
$ws.Range("E27").Value = "  +3.94%  "
$ws.Range("D28").Value = "20.24"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("D30").Value = "5.617"
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("D31").Value = "122.37"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").Value = "0.9850"
$ws.Range("E32").Value = "  +5.36%  "
$ws.Range("D33").Value = "0.09546"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").Value = "1.452"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "3.596"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "5.385"
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("D37").Value = "0.06132"
$ws.Range("E37").Value = "  +2.48%  "
$ws.Range("D38").Value = "0.02258"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").Value = "1.234"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").Value = "8.175"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").Value = "0.6047"
$ws.Range("E41").Value = "  +2.40%  "

# ---- Part 2: insert new "Frax" row at position 42, shifting rows 42-51 down to 43-52 ----
$ws.Rows.Item(42).Insert()
$ws.Range("A41").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Part 3: write final values for rows 42-51 ----
$dTextRows2 = @(42, 43, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $dTextRows2) { $ws.Cells.Item($r, 4).NumberFormat = "@" }

$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1903"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "10.32"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5726"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.260"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("A47").Value = 45
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "12.28"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.950"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "3.384"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06908"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "114.10"
$ws.Range("E51").Value = "  +5.00%  "

# ---- Part 4: remove the row that fell off the bottom (old BabyDogeCoin, now at row 52) ----
$ws.Rows.Item(52).Delete()
